$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the ImageName column (D) for each kanji row with its image filename.
$ws.Range("D3").Value  = "n5_one.jpg"
$ws.Range("D4").Value  = "n5_country.jpg"
$ws.Range("D5").Value  = "n5_person.webp"
$ws.Range("D6").Value  = "n5_year.jpg"
$ws.Range("D7").Value  = "n5_big.webp"
$ws.Range("D8").Value  = "n5_ten.jpg"
$ws.Range("D9").Value  = "n5_two.png"
$ws.Range("D10").Value = "n5_book.png"
$ws.Range("D11").Value = "n5_inside.webp"

# D9 ("two") gets a wrap-text style applied.
$ws.Range("D9").WrapText = $true

# Update the active selection.
$ws.Range("E12").Select()
